$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "target" header label to O1 (spectral clustering labels column)
$ws.Range("O1").Value = "target"

$data = New-Object 'object[,]' 17,15
$data[0,0] = 32
$data[0,1] = -0.7170444859158436
$data[0,2] = -0.9626627575504668
$data[0,3] = -1.326370565969052
$data[0,4] = -0.3508826171547983
$data[0,5] = -0.1116572877339215
$data[0,6] = 0.3025477387642772
$data[0,7] = 0.02227960656881003
$data[0,8] = -0.8292441429574355
$data[0,9] = -0.3231513895542691
$data[0,10] = -0.6051667266971232
$data[0,11] = 0.4240250088644789
$data[0,12] = 0.5689230922157072
$data[0,13] = 0.5941888100937924
$data[0,14] = 1
$data[1,0] = 21
$data[1,1] = 1.088576992162783
$data[1,2] = 1.113879672048004
$data[1,3] = 1.391852302160855
$data[1,4] = 1.097675066956605
$data[1,5] = 0.6193999948385083
$data[1,6] = 1.260478162207732
$data[1,7] = 0.2829759842181138
$data[1,8] = -0.2048732487660765
$data[1,9] = 1.280068775920671
$data[1,10] = 0.4598623887690349
$data[1,11] = 0.7524404323814086
$data[1,12] = 0.802791065925304
$data[1,13] = 0.7880756782855436
$data[1,14] = 1
$data[2,0] = 25
$data[2,1] = -1.13941715647749
$data[2,2] = -1.061880112579847
$data[2,3] = -0.3185430159423427
$data[2,4] = -0.00493742831527194
$data[2,5] = -0.3473981302948528
$data[2,6] = 0.2431232885674622
$data[2,7] = -0.497963790472687
$data[2,8] = -0.6528323163413656
$data[2,9] = -0.3835932772458677
$data[2,10] = -0.06531649189307247
$data[2,11] = 0.263425614494407
$data[2,12] = 0.3977142755814441
$data[2,13] = 0.5535369778544057
$data[2,14] = 1
$data[3,0] = 3
$data[3,1] = 1.168895327646786
$data[3,2] = 1.124790058217134
$data[3,3] = -0.108837563652328
$data[3,4] = -0.527461071075682
$data[3,5] = -0.8652942322810397
$data[3,6] = 0.1168204941201648
$data[3,7] = -0.723641543470295
$data[3,8] = -0.498731630075758
$data[3,9] = -1.025476842083122
$data[3,10] = -1.622354741293647
$data[3,11] = -2.117899344682976
$data[3,12] = -2.105917768685274
$data[3,13] = -2.604525295334807
$data[3,14] = 1
$data[4,0] = 2
$data[4,1] = 1.084148553489448
$data[4,2] = 1.096588651912747
$data[4,3] = 0.4103411666276454
$data[4,4] = 0.02692079055893389
$data[4,5] = -0.3654507667787238
$data[4,6] = 0.255985997385492
$data[4,7] = -0.170349196938015
$data[4,8] = 0.1184074125783467
$data[4,9] = -0.09407802617016045
$data[4,10] = -0.7608731938203582
$data[4,11] = -1.741872353831748
$data[4,12] = -1.936557524428362
$data[4,13] = -2.001350799166149
$data[4,14] = 1
$data[5,0] = 28
$data[5,1] = -0.8903294304691076
$data[5,2] = -1.00739468834947
$data[5,3] = -0.7484213305595308
$data[5,4] = 0.1115915246707709
$data[5,5] = 0.1855962737291569
$data[5,6] = 0.6953621033044928
$data[5,7] = 0.1702216613658133
$data[5,8] = -0.2583615203372818
$data[5,9] = 0.04432281253928385
$data[5,10] = -0.2161966413060894
$data[5,11] = 0.3946079152495584
$data[5,12] = 0.5572839498761287
$data[5,13] = 0.6264865797433866
$data[5,14] = 1
$data[6,0] = 33
$data[6,1] = -1.112188620417129
$data[6,2] = -1.052639041820789
$data[6,3] = 0.1521366645032375
$data[6,4] = -0.1063075991169151
$data[6,5] = -0.633907517591107
$data[6,6] = -0.3105400581278053
$data[6,7] = -1.076662870880767
$data[6,8] = -2.234195617877839
$data[6,9] = -2.024776068498603
$data[6,10] = -1.05364572866942
$data[6,11] = -0.397067076086726
$data[6,12] = -0.1012654246438619
$data[6,13] = -0.1105577153171547
$data[6,14] = 1
$data[7,0] = 5
$data[7,1] = 0.09099733199999999
$data[7,2] = 0.165123739
$data[7,3] = 1.033777922
$data[7,4] = 0.62781765
$data[7,5] = 0.7088679729999999
$data[7,6] = 0.9639872290000001
$data[7,7] = -0.570432089
$data[7,8] = -0.6886823870000001
$data[7,9] = -0.226930379
$data[7,10] = 0.2216372
$data[7,11] = -0.154145444
$data[7,12] = 0.244069248
$data[7,13] = 0.37543717
$data[7,14] = 1
$data[8,0] = 7
$data[8,1] = 1.109438579590896
$data[8,2] = 1.099713850078085
$data[8,3] = -0.3721542361956252
$data[8,4] = -0.8638457975380623
$data[8,5] = -1.116670404868982
$data[8,6] = -0.1384535856044998
$data[8,7] = 0.08506005601941405
$data[8,8] = -0.1424232843816636
$data[8,9] = -0.6322272669356886
$data[8,10] = -1.594814530196119
$data[8,11] = -2.152477303639709
$data[8,12] = -2.090516287548417
$data[8,13] = -2.176585305596301
$data[8,14] = 1
$data[9,0] = 16
$data[9,1] = -1.044096069467892
$data[9,2] = -0.9944797556824532
$data[9,3] = -0.4346623413973897
$data[9,4] = -0.4656609981401171
$data[9,5] = -0.3972448757409678
$data[9,6] = -0.9173023484438549
$data[9,7] = -0.6683828193368909
$data[9,8] = -0.6552500014310596
$data[9,9] = -0.7334345699974137
$data[9,10] = -1.761061959828618
$data[9,11] = -0.4595386133116331
$data[9,12] = -0.4046128729647322
$data[9,13] = 0.06517742000869589
$data[9,14] = 1
$data[10,0] = 43
$data[10,1] = -0.3543593665836199
$data[10,2] = -0.607770836055068
$data[10,3] = -1.023579362476468
$data[10,4] = -0.955669614345772
$data[10,5] = -1.000502683286009
$data[10,6] = -1.440570144316475
$data[10,7] = -0.6574608486657882
$data[10,8] = -0.3437044833041243
$data[10,9] = -0.3136140931433612
$data[10,10] = 0.1547769421110219
$data[10,11] = 0.2203860173098343
$data[10,12] = 0.1150587133144156
$data[10,13] = -0.0393656080207815
$data[10,14] = 2
$data[11,0] = 45
$data[11,1] = -0.2139746633820842
$data[11,2] = -0.7044626671279796
$data[11,3] = -2.085647036919111
$data[11,4] = -0.774573281969745
$data[11,5] = -1.204105165622195
$data[11,6] = -1.790587872433676
$data[11,7] = -0.8412445560158462
$data[11,8] = -1.00745715188306
$data[11,9] = -0.6580804567079027
$data[11,10] = -0.4286910453839183
$data[11,11] = -0.2915545232227261
$data[11,12] = -0.1628719732754185
$data[11,13] = -0.0334356291279433
$data[11,14] = 2
$data[12,0] = 13
$data[12,1] = -0.937696286
$data[12,2] = -1.015591623
$data[12,3] = -0.608602128
$data[12,4] = -1.158174718
$data[12,5] = -0.8048887790000001
$data[12,6] = -1.24772037
$data[12,7] = -1.271167965
$data[12,8] = -0.246935819
$data[12,9] = -0.041770972
$data[12,10] = -0.280753298
$data[12,11] = -0.46378129
$data[12,12] = -0.296611509
$data[12,13] = -0.032517665
$data[12,14] = 2
$data[13,0] = 55
$data[13,1] = -0.7261516585242036
$data[13,2] = -0.708529406368131
$data[13,3] = -0.4361873923187662
$data[13,4] = -0.5217290491349929
$data[13,5] = -0.6855108296681408
$data[13,6] = -1.479215207788048
$data[13,7] = -1.652506564090377
$data[13,8] = -0.3700793103049028
$data[13,9] = -0.395443007300101
$data[13,10] = -0.3867351735737296
$data[13,11] = -0.005570281690596239
$data[13,12] = 0.1983542129563711
$data[13,13] = 0.454038615761895
$data[13,14] = 2
$data[14,0] = 16
$data[14,1] = -0.898339942
$data[14,2] = -0.9166695420000001
$data[14,3] = -0.124501255
$data[14,4] = 0.389349156
$data[14,5] = 1.017939856
$data[14,6] = 1.065573911
$data[14,7] = 1.109928244
$data[14,8] = 0.875583775
$data[14,9] = 0.875780095
$data[14,10] = 1.330532029
$data[14,11] = 1.080667858
$data[14,12] = 0.87872014
$data[14,13] = 0.4082125
$data[14,14] = 3
$data[15,0] = 58
$data[15,1] = 1.112077923941515
$data[15,2] = 1.137187283022762
$data[15,3] = 1.816349257204819
$data[15,4] = 1.83863370029257
$data[15,5] = 2.121151532457759
$data[15,6] = 1.873380616886137
$data[15,7] = 2.251182816728566
$data[15,8] = 2.240347758734272
$data[15,9] = 2.209619352175399
$data[15,10] = 2.467857162570352
$data[15,11] = 1.882630143514583
$data[15,12] = 1.560192797568996
$data[15,13] = 1.36944824066759
$data[15,14] = 3
$data[16,0] = 64
$data[16,1] = 1.09846284716018
$data[16,2] = 1.123561183674226
$data[16,3] = 1.525158164547952
$data[16,4] = 1.190229985899966
$data[16,5] = 1.144727242730383
$data[16,6] = 0.9584396757252992
$data[16,7] = 0.6515114401502651
$data[16,8] = 0.6428933796817899
$data[16,9] = 0.6004911621258092
$data[16,10] = 0.9391544659798
$data[16,11] = 1.168192795304359
$data[16,12] = 1.099420116429678
$data[16,13] = 1.050510024761384
$data[16,14] = 3
$ws.Range("A2:O18").Value = $data
